# Auto-generated edit script
# Applies per-cell numeric value updates to match the target commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H17").Value = 835.1372699999999
$ws.Range("J17").Value = 835.1372699999999
$ws.Range("L17").Value = 2505.41181
$ws.Range("N17").Value = -2841.41181
$ws.Range("H33").Value = 196.94118
$ws.Range("I33").Value = 123.84615
$ws.Range("J33").Value = 434.5
$ws.Range("K33").Value = 123.84615
$ws.Range("L33").Value = 434.5
$ws.Range("M33").Value = 105.15385
$ws.Range("N33").Value = -892.5
$ws.Range("H99").Value = 817.7059
$ws.Range("I99").Value = 289.33334
$ws.Range("J99").Value = 1412.125
$ws.Range("K99").Value = 868.0000200000001
$ws.Range("L99").Value = 4236.375
$ws.Range("M99").Value = 629.9999799999999
$ws.Range("N99").Value = -7232.375
$ws.Range("H100").Value = 2233.5334
$ws.Range("I100").Value = 2240.963
$ws.Range("K100").Value = 2240.963
$ws.Range("M100").Value = -1699.963
$ws.Range("H107").Value = 2357.1904
$ws.Range("I107").Value = 3292.2856
$ws.Range("J107").Value = 1889.6428
$ws.Range("K107").Value = 3292.2856
$ws.Range("L107").Value = 1889.6428
$ws.Range("M107").Value = -1372.2856
$ws.Range("N107").Value = -5729.6428
$ws.Range("H129").Value = 655.0513
$ws.Range("I129").Value = 342.06668
$ws.Range("J129").Value = 850.6667
$ws.Range("K129").Value = 1026.20004
$ws.Range("L129").Value = 2552.0001
$ws.Range("M129").Value = 3973.79996
$ws.Range("N129").Value = -12552.0001
$ws.Range("H137").Value = 1305.2094
$ws.Range("I137").Value = 832.5
$ws.Range("K137").Value = 2497.5
$ws.Range("M137").Value = 52.5
$ws.Range("H138").Value = 1221.1464
$ws.Range("I138").Value = 804.51166
$ws.Range("J138").Value = 1680.5128
$ws.Range("K138").Value = 2413.53498
$ws.Range("L138").Value = 5041.538399999999
$ws.Range("M138").Value = 2726.46502
$ws.Range("N138").Value = -15321.5384
$ws.Range("H141").Value = 603.25
$ws.Range("I141").Value = 603.25
$ws.Range("K141").Value = 1809.75
$ws.Range("M141").Value = 3370.25

$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 4716.3335
$ws.Range("I32").Value = 4244.0166
$ws.Range("K32").Value = 4244.0166
$ws.Range("M32").Value = -3957.0166
$ws.Range("H74").Value = 1575.45
$ws.Range("I74").Value = 1093.4375
$ws.Range("J74").Value = 3503.5
$ws.Range("K74").Value = 1093.4375
$ws.Range("L74").Value = 3503.5
$ws.Range("M74").Value = -219.4375
$ws.Range("N74").Value = -5251.5
$ws.Range("H77").Value = 1575.45
$ws.Range("I77").Value = 1093.4375
$ws.Range("J77").Value = 3503.5
$ws.Range("K77").Value = 5467.1875
$ws.Range("L77").Value = 17517.5
$ws.Range("M77").Value = -1099.1875
$ws.Range("N77").Value = -26253.5
$ws.Range("H107").Value = 30000
$ws.Range("J107").Value = 30000
$ws.Range("L107").Value = 30000
$ws.Range("N107").Value = -37680
$ws.Range("H110").Value = 1477
$ws.Range("I110").Value = 836.44446
$ws.Range("J110").Value = 2197.625
$ws.Range("K110").Value = 836.44446
$ws.Range("L110").Value = 2197.625
$ws.Range("M110").Value = 1208.55554
$ws.Range("N110").Value = -6287.625
$ws.Range("H122").Value = 1564.5834
$ws.Range("I122").Value = 1724.5555
$ws.Range("J122").Value = 1084.6666
$ws.Range("K122").Value = 5173.666499999999
$ws.Range("L122").Value = 3253.9998
$ws.Range("M122").Value = -2723.666499999999
$ws.Range("N122").Value = -8153.9998
$ws.Range("H132").Value = 1126.0566
$ws.Range("I132").Value = 1054.8158
$ws.Range("K132").Value = 3164.4474
$ws.Range("M132").Value = -634.4474

$ws = $wb.Worksheets("BSM")
$ws.Range("H42").Value = 241842
$ws.Range("J42").Value = 241842
$ws.Range("L42").Value = 241842
$ws.Range("N42").Value = -242498
$ws.Range("H87").Value = 65000
$ws.Range("J87").Value = 65000
$ws.Range("L87").Value = 65000
$ws.Range("N87").Value = -67496
$ws.Range("H90").Value = 65000
$ws.Range("J90").Value = 65000
$ws.Range("L90").Value = 195000
$ws.Range("N90").Value = -207480
$ws.Range("H94").Value = 17857742
$ws.Range("J94").Value = 2502.5
$ws.Range("L94").Value = 2502.5
$ws.Range("N94").Value = -3404.5
$ws.Range("H105").Value = 71430940
$ws.Range("I105").Value = 90911290
$ws.Range("K105").Value = 90911290
$ws.Range("M105").Value = -90909543
$ws.Range("H112").Value = 33333
$ws.Range("J112").Value = 33333
$ws.Range("L112").Value = 33333
$ws.Range("N112").Value = -36287
$ws.Range("H140").Value = 20621.334
$ws.Range("J140").Value = 20621.334
$ws.Range("L140").Value = 20621.334
$ws.Range("N140").Value = -30981.334

$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 1419.9508
$ws.Range("I31").Value = 1292.2678
$ws.Range("K31").Value = 1292.2678
$ws.Range("M31").Value = -997.2678000000001
$ws.Range("H34").Value = 1419.9508
$ws.Range("I34").Value = 1292.2678
$ws.Range("K34").Value = 1292.2678
$ws.Range("M34").Value = -1090.2678
$ws.Range("H94").Value = 1283.6666
$ws.Range("I94").Value = 1452.75
$ws.Range("J94").Value = 1199.125
$ws.Range("K94").Value = 1452.75
$ws.Range("L94").Value = 1199.125
$ws.Range("M94").Value = -1001.75
$ws.Range("N94").Value = -2101.125
$ws.Range("H99").Value = 2735
$ws.Range("I99").Value = 2735
$ws.Range("K99").Value = 2735
$ws.Range("M99").Value = -1237
$ws.Range("H107").Value = 684.6
$ws.Range("I107").Value = 552.1111
$ws.Range("K107").Value = 552.1111
$ws.Range("M107").Value = 1367.8889
$ws.Range("H126").Value = 2735
$ws.Range("I126").Value = 2735
$ws.Range("K126").Value = 8205
$ws.Range("M126").Value = -5735
$ws.Range("H134").Value = 10639590
$ws.Range("I134").Value = 1361.8975
$ws.Range("K134").Value = 4085.6925
$ws.Range("M134").Value = -1550.6925
$ws.Range("H141").Value = 29576.615
$ws.Range("J141").Value = 29576.615
$ws.Range("L141").Value = 29576.615
$ws.Range("N141").Value = -39936.61500000001

$ws = $wb.Worksheets("CUL")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H25").Value = 2000
$ws.Range("I25").Value = 2000
$ws.Range("K25").Value = 6000
$ws.Range("M25").Value = -5831
$ws.Range("H30").Value = 2000
$ws.Range("I30").Value = 2000
$ws.Range("K30").Value = 6000
$ws.Range("M30").Value = -5898
$ws.Range("J131").Value = 1861.1305
$ws.Range("L131").Value = 5583.3915
$ws.Range("N131").Value = -15663.3915

$ws = $wb.Worksheets("GSM")
$ws.Range("H80").Value = 3066.6667
$ws.Range("I80").Value = 1740
$ws.Range("J80").Value = 4014.2856
$ws.Range("K80").Value = 1740
$ws.Range("L80").Value = 4014.2856
$ws.Range("M80").Value = -742
$ws.Range("N80").Value = -6010.2856
$ws.Range("H83").Value = 3066.6667
$ws.Range("I83").Value = 1740
$ws.Range("J83").Value = 4014.2856
$ws.Range("K83").Value = 8700
$ws.Range("L83").Value = 20071.428
$ws.Range("M83").Value = -3708
$ws.Range("N83").Value = -30055.428
$ws.Range("H126").Value = 2707.5
$ws.Range("I126").Value = 1812.1111
$ws.Range("J126").Value = 3858.7144
$ws.Range("K126").Value = 5436.3333
$ws.Range("L126").Value = 11576.1432
$ws.Range("M126").Value = -2966.3333
$ws.Range("N126").Value = -16516.1432
$ws.Range("H132").Value = 1638
$ws.Range("I132").Value = 1348.7407
$ws.Range("K132").Value = 4046.2221
$ws.Range("M132").Value = -1516.2221

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 2166.5
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1888
$ws.Range("H40").Value = 3166.3333
$ws.Range("I40").Value = 2799.6
$ws.Range("K40").Value = 2799.6
$ws.Range("M40").Value = -2663.6
$ws.Range("H100").Value = 1666.5555
$ws.Range("I100").Value = 1428.4286
$ws.Range("K100").Value = 1428.4286
$ws.Range("M100").Value = -887.4286
$ws.Range("H126").Value = 2166.5
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("H130").Value = 36250
$ws.Range("J130").Value = 36250
$ws.Range("L130").Value = 36250
$ws.Range("N130").Value = -46290
$ws.Range("H132").Value = 25088.023
$ws.Range("I132").Value = 1571.24
$ws.Range("J132").Value = 57750.223
$ws.Range("K132").Value = 4713.72
$ws.Range("L132").Value = 173250.669
$ws.Range("M132").Value = -2183.72
$ws.Range("N132").Value = -178310.669
$ws.Range("H136").Value = 3783.282
$ws.Range("I136").Value = 3854.111
$ws.Range("J136").Value = 2933.3333
$ws.Range("K136").Value = 11562.333
$ws.Range("L136").Value = 8799.999899999999
$ws.Range("M136").Value = -9012.332999999999
$ws.Range("N136").Value = -13899.9999

$ws = $wb.Worksheets("WVR")
$ws.Range("H132").Value = 2283.8948
$ws.Range("I132").Value = 3635.5715
$ws.Range("J132").Value = 1495.4166
$ws.Range("K132").Value = 10906.7145
$ws.Range("L132").Value = 4486.2498
$ws.Range("M132").Value = -8376.7145
$ws.Range("N132").Value = -9546.2498
$ws.Range("H136").Value = 540.9032
$ws.Range("I136").Value = 496.58334
$ws.Range("K136").Value = 1489.75002
$ws.Range("M136").Value = 1060.24998

